# Append 6 new instrument rows (rows 5-10) to the "Instrumentos" sheet,
# matching columns: nombre | categoria | estado | ubicacion

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Vidrios de reloj", "Química", "DISPONIBLE", "Laboratorio 10"),
    @("Desfibrilador externo automático", "Biomédica", "MANTENIMIENTO", "Laboratorio 31"),
    @("Portaobjetos", "Química", "DISPONIBLE", "Laboratorio 22"),
    @("Monitor de signos vitales", "Biomédica", "DISPONIBLE", "Laboratorio 31"),
    @("Osciloscopio", "Electricidad", "PRESTADO", "Laboratorio 22"),
    @("Cautín", "Electricidad", "DISPONIBLE", "Laboratorio A")
)

$row = 5
foreach ($entry in $newRows) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]
    $row++
}
